# Update the "dSF" column (F) values on Sheet1 to re-pulled / re-pushed data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value  = -4
$ws.Range("F4").Value  = -7
$ws.Range("F6").Value  = -2
$ws.Range("F7").Value  = 0
$ws.Range("F9").Value  = -3
$ws.Range("F11").Value = 1
$ws.Range("F12").Value = -2
